# #5: cash & deposit done
# Sheet "存款" (deposit): add bank/deposit_type/currency header labels and the
# shared property_category/category/date/legislator_name/legislator_id/
# source_file/index metadata columns (G:M) that the other sheets already carry.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("存款")

# ---- Header row (row 1) ----
$ws.Cells.Item(1, 2).Value = "bank"
$ws.Cells.Item(1, 3).Value = "deposit_type"
$ws.Cells.Item(1, 4).Value = "currency"
$ws.Cells.Item(1, 5).Value = "owner"
$ws.Cells.Item(1, 6).Value = "total"
$ws.Cells.Item(1, 7).Value = "property_category"
$ws.Cells.Item(1, 8).Value = "category"
$ws.Cells.Item(1, 9).Value = "date"
$ws.Cells.Item(1, 10).Value = "legislator_name"
$ws.Cells.Item(1, 11).Value = "legislator_id"
$ws.Cells.Item(1, 12).Value = "source_file"
$ws.Cells.Item(1, 13).Value = "index"

# ---- Data rows (2-16): append the constant metadata columns G:M ----
for ($row = 2; $row -le 16; $row++) {
    $idx = $ws.Cells.Item($row, 1).Value()

    $ws.Cells.Item($row, 7).Value = "deposit"
    $ws.Cells.Item($row, 8).Value = "normal"
    $ws.Cells.Item($row, 9).Value = "2013-12-12"
    $ws.Cells.Item($row, 10).Value = "林鴻池"
    $ws.Cells.Item($row, 11).Value = 1340
    $ws.Cells.Item($row, 12).Value = "tmp67ea1"
    $ws.Cells.Item($row, 13).Value = $idx
}
